# ---------------------------------------------------------------------------
# Applies the "finished second corr calcuations" edit:
#  - NYT Books sheet: change the active selection
#  - Correlation sheet: rename to "Box Office $ Correlation"
#  - Add a new "Rank + In Theaters Corr" sheet (copy of the Weekly Data
#    week/rank columns plus an "In Theaters?" flag and two CORREL() results)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. NYT Books: just move the selection -------------------------------
$nytBooks = $wb.Worksheets.Item("NYT Books")
$nytBooks.Range("B1:E21").Select() | Out-Null

# --- 2. Rename the Correlation sheet --------------------------------------
$corr = $wb.Worksheets.Item("Correlation")
$corr.Name = "Box Office `$ Correlation"

# --- 3. Add the new "Rank + In Theaters Corr" sheet as the last sheet ----
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Rank + In Theaters Corr"

# Header row
$newSheet.Range("A1").Value = "Week Start"
$newSheet.Range("B1").Value = "Week End"
$newSheet.Range("C1").Value = "Rank"
$newSheet.Range("D1").Value = "In Theaters? (1 = Y, 0 = N)"
$newSheet.Range("F1").Value = "Correlation:"

# Data rows 2-21: Week End (B), Rank (C), In Theaters? (D)
$bVals = @(42966,42973,42980,42987,42994,43001,43008,43015,43022,43029,43036,43043,43050,43057,43064,43071,43078,43085,43092,43099)
$cVals = @(6,5,4,1,1,1,1,2,2,2,3,3,3,3,3,4,4,6,5,5)
$dVals = @(0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0)

for ($i = 0; $i -lt 20; $i++) {
    $r = $i + 2
    $newSheet.Cells.Item($r, 2).Value = $bVals[$i]
    $newSheet.Cells.Item($r, 3).Value = $cVals[$i]
    $newSheet.Cells.Item($r, 4).Value = $dVals[$i]
}

# Week Start (A) = Week End (B) - 6, as a formula (shared across A3:A21)
$newSheet.Range("A2").Formula = "=B2-6"
$newSheet.Range("A3:A21").Formula = "=B3-6"

# Date formatting for the Week Start / Week End columns (A & B)
$newSheet.Range("A2:B21").NumberFormat = "yyyy\-mm\-dd;@"
$newSheet.Range("B1").NumberFormat = "yyyy\-mm\-dd;@"

# Correlation formulas
$newSheet.Range("F2").Formula = "=CORREL(C2:C21,D2:D21)"
$newSheet.Range("F4").Value = "Correlation once first movie opens:"
$newSheet.Range("F5").Formula = "=CORREL(C5:C21,D5:D21)"

# Column widths, roughly matching the source workbook
$newSheet.Columns("B:C").ColumnWidth = 10.83
$newSheet.Columns("D").ColumnWidth = 22.5

# Selection / active cell on the new sheet
$newSheet.Range("F6").Select() | Out-Null
$newSheet.Activate()
